$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.587199866771698
$ws.Range("B1").Value = 2.685224294662476
$ws.Range("C1").Value = 6.423510551452637
$ws.Range("D1").Value = 1.85191535949707
$ws.Range("E1").Value = 1.586808681488037
